$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (in-place run edits, preserve rich-text runs) ---
# A8: "Volume 32   Number  25" -> "...26"  (replace just the trailing "25")
$ws.Range("A8").Characters(21, 2).Text = "26"

# C9: "Report Covering the Week  6/16/2025  Through  6/22/2025"
#     -> week shifted forward by one week
$ws.Range("C9").Characters(27, 9).Text = "6/23/2025"
$ws.Range("C9").Characters(47, 9).Text = "6/29/2025"

# --- Numeric cell updates (rows 15-30 weekly crime-stat tables) ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 19
$ws.Range("K15").Value = -21.052631578947
$ws.Range("L15").Value = 15.384615384615
$ws.Range("N15").Value = -6.25
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -62.5
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = 6.25
$ws.Range("I16").Value = 120
$ws.Range("J16").Value = 127
$ws.Range("K16").Value = -5.511811023622
$ws.Range("L16").Value = -5.511811023622
$ws.Range("M16").Value = -7.692307692307
$ws.Range("N16").Value = -65.417867435158
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -30.769230769230
$ws.Range("F17").Value = 51
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = 30.769230769230
$ws.Range("I17").Value = 254
$ws.Range("J17").Value = 222
$ws.Range("K17").Value = 14.414414414414
$ws.Range("L17").Value = 22.115384615384
$ws.Range("M17").Value = 109.917355371901
$ws.Range("N17").Value = 66.013071895424
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 116.666666666667
$ws.Range("I18").Value = 122
$ws.Range("J18").Value = 87
$ws.Range("K18").Value = 40.229885057471
$ws.Range("L18").Value = 6.086956521739
$ws.Range("M18").Value = -33.695652173913
$ws.Range("N18").Value = -83.754993342210
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -26.415094339622
$ws.Range("I19").Value = 369
$ws.Range("J19").Value = 436
$ws.Range("K19").Value = -15.366972477064
$ws.Range("L19").Value = 17.891373801916
$ws.Range("M19").Value = 90.206185567010
$ws.Range("N19").Value = 27.681660899654
$ws.Range("C20").Value = 10
$ws.Range("E20").Value = -23.076923076923
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = -32
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 231
$ws.Range("K20").Value = 7.792207792207
$ws.Range("L20").Value = -6.741573033707
$ws.Range("M20").Value = 107.5
$ws.Range("N20").Value = -72.846237731733
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = -30.769230769230
$ws.Range("F21").Value = 155
$ws.Range("G21").Value = 169
$ws.Range("H21").Value = -8.284023668639
$ws.Range("I21").Value = 1131
$ws.Range("J21").Value = 1125
$ws.Range("K21").Value = 0.533333333333
$ws.Range("L21").Value = 8.126195028680
$ws.Range("M21").Value = 47.650130548302
$ws.Range("N21").Value = -54.358353510895
$ws.Range("F22").Value = 1
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 12
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 58
$ws.Range("J23").Value = 56
$ws.Range("K23").Value = 3.571428571428
$ws.Range("L23").Value = -18.309859154929
$ws.Range("M23").Value = 100
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 106.25
$ws.Range("F24").Value = 143
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 76.543209876543
$ws.Range("I24").Value = 823
$ws.Range("J24").Value = 689
$ws.Range("K24").Value = 19.448476052249
$ws.Range("L24").Value = 6.193548387096
$ws.Range("M24").Value = 92.740046838407
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 450
$ws.Range("F25").Value = 66
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 175
$ws.Range("I25").Value = 324
$ws.Range("J25").Value = 277
$ws.Range("K25").Value = 16.967509025270
$ws.Range("L25").Value = 1.886792452830
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 60
$ws.Range("G26").Value = 68
$ws.Range("H26").Value = -11.764705882352
$ws.Range("I26").Value = 357
$ws.Range("J26").Value = 293
$ws.Range("K26").Value = 21.843003412969
$ws.Range("L26").Value = 31.734317343173
$ws.Range("M26").Value = 14.790996784565
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = -18.181818181818
$ws.Range("L27").Value = -5.263157894736
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = 9.677419354838
$ws.Range("L28").Value = -17.073170731707
$ws.Range("M29").Value = -73.333333333333
$ws.Range("N29").Value = -83.333333333333
$ws.Range("M30").Value = -75
$ws.Range("N30").Value = -85

# --- Cells that flip from a numeric 0-count to the "no data" text markers ---
# These reuse the existing shared strings ("0" = index 20, "***.*" = index 21)
# that are already used elsewhere in these same rows (e.g. D22/E22, D31/E31).
# Assigning a leading apostrophe forces Excel to store the value as text
# instead of as the number 0; we then copy the number-format/style from a
# neighboring cell that already carries the correct "text label" style so
# the cell's style index matches the rest of the row.

$ws.Range("C22").Value = "'0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("C28").Value = "'0"
$ws.Range("M28").Copy()
$ws.Range("C28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("G31").Value = "'0"
$ws.Range("F31").Copy()
$ws.Range("G31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("H31").Value = "***.*"
$ws.Range("F31").Copy()
$ws.Range("H31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$excel.CutCopyMode = 0
